$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Date" column (C) holds values that look like dates ("12/30/2024"),
# but the source data keeps them as plain text. Force text formatting on
# the range before writing so the engine doesn't coerce them into date
# serial numbers, then reset the style afterwards so no visible style
# index is left behind on the cells.
$dateRange = $ws.Range("C24:C33")
$dateRange.NumberFormat = "@"

# Row 24 used to be the last original row (order 53d29848..., "Chironji
# Seeds (4 x 500gm)", dated 12/31/2024). It is now overwritten by a new
# order placed on 12/30/2024.
$ws.Range("A24").Value = "f4997297-362d-450a-b1d9-320725609ea5"
$ws.Range("B24").Value = "Fresh and Natural Black Seeds (5 x 1kg)"
$ws.Range("C24").Value = "12/30/2024"
$ws.Range("D24").Value = 1710
$ws.Range("E24").Value = "success"

# New rows 25-32: newly checked-out orders, all dated 12/30/2024.
$ws.Range("A25").Value = "791c62fd-46bb-46e4-9d54-7239b14c263e"
$ws.Range("B25").Value = "Fresh and Natural Black Seeds (2 x 1kg)"
$ws.Range("C25").Value = "12/30/2024"
$ws.Range("D25").Value = 684
$ws.Range("E25").Value = "success"

$ws.Range("A26").Value = "742fdec7-21fd-4e77-b46c-87fa033b9fd1"
$ws.Range("B26").Value = "Fresh and Natural Black Seeds (5 x 250gm)"
$ws.Range("C26").Value = "12/30/2024"
$ws.Range("D26").Value = 3758
$ws.Range("E26").Value = "success"

$ws.Range("A27").Value = "742fdec7-21fd-4e77-b46c-87fa033b9fd1"
$ws.Range("B27").Value = "Chironji Seeds (5 x 250gm)"
$ws.Range("C27").Value = "12/30/2024"
$ws.Range("D27").Value = 3758
$ws.Range("E27").Value = "success"

$ws.Range("A28").Value = "742fdec7-21fd-4e77-b46c-87fa033b9fd1"
$ws.Range("B28").Value = "Chironji Seeds (1 x 1kg)"
$ws.Range("C28").Value = "12/30/2024"
$ws.Range("D28").Value = 3758
$ws.Range("E28").Value = "success"

$ws.Range("A29").Value = "6b4532c9-bd83-460d-803a-410e4deddf68"
$ws.Range("B29").Value = "Fresh and Natural Black Seeds (5 x 250gm)"
$ws.Range("C29").Value = "12/30/2024"
$ws.Range("D29").Value = 5596
$ws.Range("E29").Value = "success"

$ws.Range("A30").Value = "6b4532c9-bd83-460d-803a-410e4deddf68"
$ws.Range("B30").Value = "Chironji Seeds (5 x 250gm)"
$ws.Range("C30").Value = "12/30/2024"
$ws.Range("D30").Value = 5596
$ws.Range("E30").Value = "success"

$ws.Range("A31").Value = "6b4532c9-bd83-460d-803a-410e4deddf68"
$ws.Range("B31").Value = "Chironji Seeds (2 x 1kg)"
$ws.Range("C31").Value = "12/30/2024"
$ws.Range("D31").Value = 5596
$ws.Range("E31").Value = "success"

$ws.Range("A32").Value = "2c219e31-f87a-4a04-9443-2feaa9f93d3b"
$ws.Range("B32").Value = "Chironji Seeds (3 x 1kg)"
$ws.Range("C32").Value = "12/30/2024"
$ws.Range("D32").Value = 4764
$ws.Range("E32").Value = "success"

# Row 33 now holds the order that used to sit on row 24 (last original
# row), dated 12/31/2024.
$ws.Range("A33").Value = "53d29848-bbe1-4354-9ef1-50a1cdb7ad7b"
$ws.Range("B33").Value = "Chironji Seeds (4 x 500gm)"
$ws.Range("C33").Value = "12/31/2024"
$ws.Range("D33").Value = 2985.44
$ws.Range("E33").Value = "success"

# Drop the text-number-format override now that all the date-looking
# strings are safely stored as text, so no stray style index lingers on
# the cells.
$dateRange.Style = "Normal"

# The old Grand Total row (26) moves down to row 35 (row 34 stays blank,
# and row 26's old "Grand Total" text is naturally overwritten above by
# the new order data), and the total itself increases to reflect the new
# orders.
$ws.Range("D35").Value = "Grand Total: ₹53181.94"
